$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "Jorge Aquino-Motores de aplicação" class from the 9:50/10:40
# slots (rows 6 & 7) to the 7:50/8:40 slots (rows 3 & 4) on Monday column (B)
$ws.Range("B3").Value = "Jorge Aquino-Motores de aplicação"
$ws.Range("B4").Value = "Jorge Aquino-Motores de aplicação"
$ws.Range("B6").Value = "-"
$ws.Range("B7").Value = "-"
